$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task description labels (column A) ---
# Rows 3-8 keep (or restate) their existing labels; rows 9-22 get new/updated
# task descriptions reflecting the finer-grained sprint 4 breakdown.
$ws.Range("A3").Value = 'Update Noteviewer Desktop Saul'
$ws.Range("A4").Value = 'Update Noteviewer Web Saul'
$ws.Range("A5").Value = 'Fix Expanded Note Desktop Corey'
$ws.Range("A6").Value = 'Fix Expanded Note Web Corey'
$ws.Range("A7").Value = 'Update UI Desktop Saul'
$ws.Range("A8").Value = 'Update UI Web Corey'
$ws.Range("A9").Value = 'Export a list of sources Web Corey'
$ws.Range("A10").Value = 'Export a list of sources Desktop Corey'
$ws.Range("A11").Value = 'create basic layout Project Page Web Saul'
$ws.Range("A12").Value = 'create basic layout Project Page Desktop Saul'
$ws.Range("A13").Value = 'implement add source project page web saul'
$ws.Range("A14").Value = 'implement add source project page desktop saul'
$ws.Range("A15").Value = 'implement remove source project page web saul'
$ws.Range("A16").Value = 'implement remove source project page desktop saul'
$ws.Range("A17").Value = 'create export list of sources dialog web corey'
$ws.Range("A18").Value = 'create export list of sources dialog desktop corey'
$ws.Range("A19").Value = 'implment the formated sources web corey '
$ws.Range("A20").Value = 'implment the formated sources desktop corey '
$ws.Range("A21").Value = 'update class diagram for web corey'
$ws.Range("A22").Value = 'fix sprint backlog 4 diagram corey'

# --- Burndown values: rows 3-20 are fully re-estimated at 3 across all 14 days ---
$ws.Range("C3:Q20").Value = 3

# --- Rows 21-22 (new tasks) ramp down to 0 starting day 8 (column J) ---
$ws.Range("C21:I21").Value = 2
$ws.Range("J21:Q21").Value = 0
$ws.Range("C22:I22").Value = 1
$ws.Range("J22:Q22").Value = 0

# --- View: zoom level and active selection ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("S20").Select()
